$wb = $excel.ActiveWorkbook

# This script applies a scheduled market-data refresh to the per-job "Profits" sheets.
# Each sheet (one per Disciple of the Hand job) tracks current market-board prices
# (columns H-L) and the resulting Leve profit margins (columns M/N). Updated rows below.

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2243.4055
$ws.Range("J17").Value = 2243.4055
$ws.Range("L17").Value = 6730.2165
$ws.Range("N17").Value = -7066.2165
$ws.Range("H19").Value = 1445.625
$ws.Range("I19").Value = 1489
$ws.Range("J19").Value = 1419.6
$ws.Range("K19").Value = 1489
$ws.Range("L19").Value = 1419.6
$ws.Range("M19").Value = -1314
$ws.Range("N19").Value = -1769.6
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 100000
$ws.Range("M65").Value = -96880
$ws.Range("H74").Value = 4939.75
$ws.Range("I74").Value = 4939.75
$ws.Range("K74").Value = 4939.75
$ws.Range("M74").Value = -4003.75
$ws.Range("H77").Value = 4939.75
$ws.Range("I77").Value = 4939.75
$ws.Range("K77").Value = 24698.75
$ws.Range("M77").Value = -20018.75
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("H111").Value = 415.6
$ws.Range("I111").Value = 374
$ws.Range("J111").Value = 443.33334
$ws.Range("K111").Value = 1122
$ws.Range("L111").Value = 1330.00002
$ws.Range("M111").Value = 1945
$ws.Range("N111").Value = -7464.000019999999
$ws.Range("H112").Value = 2558.7896
$ws.Range("J112").Value = 2565.4443
$ws.Range("L112").Value = 7696.3329
$ws.Range("N112").Value = -9912.332900000001
$ws.Range("H116").Value = 83111.664
$ws.Range("I116").Value = 78500.625
$ws.Range("J116").Value = 120000
$ws.Range("K116").Value = 78500.625
$ws.Range("L116").Value = 120000
$ws.Range("M116").Value = -75058.625
$ws.Range("N116").Value = -126884
$ws.Range("H132").Value = 4986.875
$ws.Range("I132").Value = 4652.4
$ws.Range("K132").Value = 13957.2
$ws.Range("M132").Value = -11427.2
$ws.Range("H135").Value = 2176
$ws.Range("I135").Value = 354.14285
$ws.Range("J135").Value = 14929
$ws.Range("K135").Value = 3187.28565
$ws.Range("L135").Value = 134361
$ws.Range("M135").Value = -652.2856500000003
$ws.Range("N135").Value = -139431
$ws.Range("H137").Value = 2280.6365
$ws.Range("I137").Value = 2662.4285
$ws.Range("K137").Value = 7987.2855
$ws.Range("M137").Value = -5437.2855
$ws.Range("H141").Value = 7998.5
$ws.Range("I141").Value = 7998.5
$ws.Range("K141").Value = 23995.5
$ws.Range("M141").Value = -18815.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9879.200000000001
$ws.Range("I2").Value = 9865.666999999999
$ws.Range("K2").Value = 9865.666999999999
$ws.Range("M2").Value = -9752.666999999999
$ws.Range("H61").Value = 4618.778
$ws.Range("I61").Value = 5010
$ws.Range("J61").Value = 3249.5
$ws.Range("K61").Value = 5010
$ws.Range("L61").Value = 3249.5
$ws.Range("M61").Value = -4798
$ws.Range("N61").Value = -3673.5
$ws.Range("H74").Value = 5764.8687
$ws.Range("I74").Value = 6760.4136
$ws.Range("J74").Value = 2557
$ws.Range("K74").Value = 6760.4136
$ws.Range("L74").Value = 2557
$ws.Range("M74").Value = -5886.4136
$ws.Range("N74").Value = -4305
$ws.Range("H77").Value = 5764.8687
$ws.Range("I77").Value = 6760.4136
$ws.Range("J77").Value = 2557
$ws.Range("K77").Value = 33802.068
$ws.Range("L77").Value = 12785
$ws.Range("M77").Value = -29434.068
$ws.Range("N77").Value = -21521
$ws.Range("H116").Value = 9879.200000000001
$ws.Range("I116").Value = 9865.666999999999
$ws.Range("K116").Value = 9865.666999999999
$ws.Range("M116").Value = -7571.666999999999
$ws.Range("H136").Value = 4618.778
$ws.Range("I136").Value = 5010
$ws.Range("J136").Value = 3249.5
$ws.Range("K136").Value = 15030
$ws.Range("L136").Value = 9748.5
$ws.Range("M136").Value = -12480
$ws.Range("N136").Value = -14848.5

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9879.200000000001
$ws.Range("I3").Value = 9865.666999999999
$ws.Range("K3").Value = 9865.666999999999
$ws.Range("M3").Value = -9751.666999999999
$ws.Range("H134").Value = 4584.9414
$ws.Range("I134").Value = 1917.5
$ws.Range("K134").Value = 5752.5
$ws.Range("M134").Value = -3217.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2307.8333
$ws.Range("J99").Value = 1999
$ws.Range("L99").Value = 1999
$ws.Range("N99").Value = -4995
$ws.Range("H126").Value = 2307.8333
$ws.Range("J126").Value = 1999
$ws.Range("L126").Value = 5997
$ws.Range("N126").Value = -10937
$ws.Range("H134").Value = 8050
$ws.Range("I134").Value = 8050
$ws.Range("K134").Value = 24150
$ws.Range("M134").Value = -21615

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 149.2
$ws.Range("I23").Value = 73
$ws.Range("J23").Value = 168.25
$ws.Range("K23").Value = 219
$ws.Range("L23").Value = 504.75
$ws.Range("M23").Value = 16
$ws.Range("N23").Value = -974.75
$ws.Range("H113").Value = 1525
$ws.Range("J113").Value = 1750
$ws.Range("L113").Value = 5250
$ws.Range("N113").Value = -9590

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3896.182
$ws.Range("I132").Value = 3580.1428
$ws.Range("J132").Value = 4449.25
$ws.Range("K132").Value = 10740.4284
$ws.Range("L132").Value = 13347.75
$ws.Range("M132").Value = -8210.428400000001
$ws.Range("N132").Value = -18407.75

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 9873
$ws.Range("J22").Value = 7972.1113
$ws.Range("L22").Value = 7972.1113
$ws.Range("N22").Value = -8562.1113
$ws.Range("H27").Value = 9873
$ws.Range("J27").Value = 7972.1113
$ws.Range("L27").Value = 7972.1113
$ws.Range("N27").Value = -8186.1113
$ws.Range("H55").Value = 1333
$ws.Range("I55").Value = 1999
$ws.Range("K55").Value = 1999
$ws.Range("M55").Value = -1826
$ws.Range("H132").Value = 4519.8945
$ws.Range("I132").Value = 2260.7273
$ws.Range("K132").Value = 6782.1819
$ws.Range("M132").Value = -4252.1819

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2630.1
$ws.Range("I132").Value = 1672
$ws.Range("K132").Value = 5016
$ws.Range("M132").Value = -2486
$ws.Range("H135").Value = 164949
$ws.Range("J135").Value = 164949
$ws.Range("L135").Value = 164949
$ws.Range("N135").Value = -175089
$ws.Range("H136").Value = 1879.1111
$ws.Range("I136").Value = 1935.0416
$ws.Range("J136").Value = 1431.6666
$ws.Range("K136").Value = 5805.1248
$ws.Range("L136").Value = 4294.9998
$ws.Range("M136").Value = -3255.1248
$ws.Range("N136").Value = -9394.9998
